$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of an existing header cell (bold, centered, bordered)
# so the new header cells match the look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins / Losses / Ties) for every player row.
$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 84  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 78  # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF -> Ties
}
